$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1169995834814548
$ws.Range("C2").Value = 109.9114832445916
$ws.Range("D2").Value = 19575605.8673771
$ws.Range("E2").Value = 14773364.14517103
$ws.Range("G2").Value = 34349080.04103096
